$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 36 (shifts old row 36 data down to row 38)
$ws.Rows("36:37").Insert()

# Populate the static (non-varying) columns for the two newly inserted rows
# by copying the full row from row 35, which has identical static values across all data rows
$ws.Range("A35:R35").Copy($ws.Range("A36:R36"))
$ws.Range("A35:R35").Copy($ws.Range("A37:R37"))

# Update the varying columns (D, J, K, L, M, P) for rows 13 through 38
$ws.Range("D13").Value = 44424
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 2500
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = 2750
$ws.Range("P13").Value = 1375
$ws.Range("D14").Value = 44279
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 1700
$ws.Range("L14").Value = 1800
$ws.Range("M14").Value = 1750
$ws.Range("P14").Value = 875
$ws.Range("D15").Value = 44356
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 2400
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2450
$ws.Range("P15").Value = 1225
$ws.Range("D16").Value = 44221
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 2900
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 2950
$ws.Range("P16").Value = 1475
$ws.Range("D17").Value = 44323
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = 2450
$ws.Range("P17").Value = 1225
$ws.Range("D18").Value = 44349
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 1800
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 1900
$ws.Range("P18").Value = 950
$ws.Range("D19").Value = 44258
$ws.Range("J19").Value = 150
$ws.Range("K19").Value = 2400
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = 2450
$ws.Range("P19").Value = 1225
$ws.Range("D20").Value = 44243
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 2900
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = 2950
$ws.Range("P20").Value = 1475
$ws.Range("D21").Value = 44176
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 1900
$ws.Range("L21").Value = 2000
$ws.Range("M21").Value = 1950
$ws.Range("P21").Value = 975
$ws.Range("D22").Value = 44237
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 2500
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = 2750
$ws.Range("P22").Value = 1375
$ws.Range("D23").Value = 44298
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 1400
$ws.Range("L23").Value = 1500
$ws.Range("M23").Value = 1450
$ws.Range("P23").Value = 725
$ws.Range("D24").Value = 44169
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 2500
$ws.Range("M24").Value = 2250
$ws.Range("P24").Value = 1125
$ws.Range("D25").Value = 44410
$ws.Range("J25").Value = 250
$ws.Range("K25").Value = 2800
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = 2900
$ws.Range("P25").Value = 1450
$ws.Range("D26").Value = 44320
$ws.Range("J26").Value = 250
$ws.Range("K26").Value = 1400
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = 1450
$ws.Range("P26").Value = 725
$ws.Range("D27").Value = 44405
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 3800
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = 3900
$ws.Range("P27").Value = 1950
$ws.Range("D28").Value = 44305
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 900
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 950
$ws.Range("P28").Value = 475
$ws.Range("D29").Value = 44417
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 3500
$ws.Range("M29").Value = 3250
$ws.Range("P29").Value = 1625
$ws.Range("D30").Value = 44181
$ws.Range("J30").Value = 250
$ws.Range("K30").Value = 1400
$ws.Range("L30").Value = 1500
$ws.Range("M30").Value = 1450
$ws.Range("P30").Value = 725
$ws.Range("D31").Value = 44319
$ws.Range("J31").Value = 300
$ws.Range("K31").Value = 1900
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = 1950
$ws.Range("P31").Value = 975
$ws.Range("D32").Value = 44253
$ws.Range("J32").Value = 300
$ws.Range("K32").Value = 2400
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = 2450
$ws.Range("P32").Value = 1225
$ws.Range("D33").Value = 44272
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 2800
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = 2900
$ws.Range("P33").Value = 1450
$ws.Range("D34").Value = 44370
$ws.Range("J34").Value = 400
$ws.Range("K34").Value = 3400
$ws.Range("L34").Value = 3500
$ws.Range("M34").Value = 3445
$ws.Range("P34").Value = 1722
$ws.Range("D35").Value = 44326
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 1400
$ws.Range("L35").Value = 1500
$ws.Range("M35").Value = 1450
$ws.Range("P35").Value = 725
$ws.Range("D36").Value = 44343
$ws.Range("J36").Value = 300
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = 1750
$ws.Range("P36").Value = 875
$ws.Range("D37").Value = 44284
$ws.Range("J37").Value = 300
$ws.Range("K37").Value = 1800
$ws.Range("L37").Value = 2000
$ws.Range("M37").Value = 1900
$ws.Range("P37").Value = 950
$ws.Range("D38").Value = 44334
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 2800
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 2900
$ws.Range("P38").Value = 1450
